$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (A1 = 0, bold+bordered style); this shifts the
# shared-string question-bank cell from A2 up to A1 and drops the style ref.
$ws.Rows.Item(1).Delete()

$newText = @'
questions = [
    {
        "title": "You are a Blazor developer working on a project. During the initialization phase of a Blazor component, you need to perform a specific action to set up initial state and configurations. You want to ensure that this action is only executed once, when the component is first rendered. Which action should you take?",
        "ques_type": 2,
        "options": [
            "Implement the OnInitialized lifecycle method.",
            "Implement the OnParametersSetAsync lifecycle method.",
            "Implement the OnParametersSet method.",
            "Use the OnAfterRender lifecycle method."
        ],
        "score": "Implement the OnInitialized lifecycle method."
    },
    {
        "title": "You need to set up route templates for different pages in a Blazor application. You must establish routes that match specific URL patterns and direct users to the corresponding components.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Implement route templates using inline attributes in HTML markup.",
            "Specify route templates in a separate configuration file.",
            "Define route patterns using the @page directive in Razor components.",
            "Configure route templates using a built-in routing component in Blazor."
        ],
        "score": "Define route patterns using the @page directive in Razor components."
    },
    {
        "title": "You are a Blazor developer working on a complex application with multiple interconnected components, one of which handles user authentication status. After a successful login, you need to manage the authentication state across various parts of the application.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Trigger a complete page refresh after a successful login.",
            "Use session storage.",
            "Use the local component state.",
            "Use a centralized application state."
        ],
        "score": "Use a centralized application state."
    },
    {
        "title": "You are a Blazor developer working on integrating JavaScript functionalities into your Blazor application. You need to call a JavaScript function from your Blazor component to handle a specific task. Which action should you take?",
        "ques_type": 2,
        "options": [
            "Use the InvokeAsync method.",
            "Use the ExecuteAsync method.",
            "Use the JSInterop attribute.",
            "Use the JSRuntime service."
        ],
        "score": "Use the JSRuntime service."
    }
]
'@

$ws.Range("A1").Value2 = $newText

# Embedded newlines make Excel auto-grow the row height on write; AutoFit
# re-measures it back down, which drops the ht/customHeight override so the
# row stays at the sheet's default height (matches the original formatting).
$ws.Rows.Item(1).AutoFit()
